# Update training_set_prediction sheet: reorder rows and refresh
# Target/Prediction values per "update main and new scripts" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FNQWTTWCYHHMVPYCDYCHFKR"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0

$ws.Range("A3").Value = "MEANSRVMVRVLLLALVVQVTLSQHWSYGWLPGGKRSVGELEATIRMMGTGEVVSLPEEASAQTQERLRPYNVINDDSSHFDRKKRSPNK"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0

$ws.Range("A4").Value = "MASVVPLKEKKLLEVKLGELPSWILMRDFTPSGIAGAFQRGYYRYYNKYVNVKKGSIAGLSMVLAAYVFLNYCRSYKELKHERLRKYH"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0

$ws.Range("A5").Value = "SILSGNFGVGKKIVCGLSGLC"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1

$ws.Range("A6").Value = "ASNQDFMRF"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 1

$ws.Range("A7").Value = "LFKLLGKIIHHVGNFVHGFSHVF"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 1

$ws.Range("A8").Value = "RPRPNYRPRPIYRP"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 1

$ws.Range("A9").Value = "MKTFSVAVAVAVVLAFICTQESSALPVTGIEELVEPVSSDNNDNHQGLPVELRERLVNIRKKRAPTDCIPYCYPTGDGFHCGVTCRF"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 1

$ws.Range("A10").Value = "ENKYFSQVVITTQCDEHRQLQRANVQWDEEVSQYF"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0

$ws.Range("A11").Value = "KFYFTFPS"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 1

$ws.Range("A12").Value = "GILSTFKGLAKGVAKDLAGNLLDKFKCKITGC"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 1

$ws.Range("A13").Value = "KSYGNGVHCNKKKCWVDWGSAISTIGNNSAANWATGGAAGWKS"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 1

$ws.Range("A14").Value = "MDGKAPAAFVEPGEFNEVMKRLDQIDEKVEFVNSEVAQRIGKKVGRDIGILYGGVIGLLLFLIYVQISSMFM"
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 1

$ws.Range("A15").Value = "RSALSCQMCELVVKKYEGSADKDANVIKKDFDAECKKLFHTIPFGTRECDHYVNSKVDPIIHELEGGTAPKDVCTKLNECP"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 1

$ws.Range("A16").Value = "NGMYFFYLNSIPAEMGRQCCAHADTYIYAMERVMFVPQFVCNSGIWTGHWKPAFRLPECYSTPWWKKS"
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 0

$ws.Range("A17").Value = "MHLSPQEKDKLLIVTAALLAERRLNRGLKLNHPEAVAWLSFLVLEGARDGKSVAELMQEGTTWLSRNQVMDGIPELVQEVQIEAVFPDGTKLVTLHDPIR"
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 0

$ws.Range("A18").Value = "WFDVDLNNIQGWIAITDGLFLEEYNKACWSCQGGPQTIHMCIHDVLIHQPFTPHAAL"
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 0

$ws.Range("A19").Value = "GGYYCPFRQDKCHRHCRSFGRKAGYCGGFLKKTCICV"
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = 1

$ws.Range("A20").Value = "MPKLAVVLLVLLILPLSYFDAAGGQAVQWDRRGNGLARYLQRGDRDVRECQVDTPGSSWGKCCMTRMCGTMCCSRSVCTCVYHWRRGHGCSCPG"
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 0

$ws.Range("A21").Value = "MQQEALGMVETKGLTAAIEAADAMVKSANVMLVGYEKIGSGLVTVIVRGDVGAVKAATDAGAAAARNVGEVKAVHVIPRPHTDVEKILPKGISQ"
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 1

$ws.Range("A22").Value = "QQDYTGWFDF"
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 1

$ws.Range("A23").Value = "GLFNVFKGLKTAGKHVAGSLLNQLKCKVSGGC"
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 1

$ws.Range("A24").Value = "LICVKEKFLFSETTETCPDGQNVCFNQAHLIYPGKYKRTRGCAATCPKLQNRDVIFCCSTDKCNL"
$ws.Range("B24").Value = 0
$ws.Range("C24").Value = 1

$ws.Range("A25").Value = "FLPAIAGMAAKFLPKIFCAISKKC"
$ws.Range("B25").Value = 1
$ws.Range("C25").Value = 1

$ws.Range("A26").Value = "SMWSGMWRRKLKKLRNALKKKLKGEK"
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = 1

$ws.Range("A27").Value = "VTSYTLSDVVPLKDVVPEWVRIGFSATPGAEYAAHEVLSWSFHSELSGTSSKQ"
$ws.Range("B27").Value = 0
$ws.Range("C27").Value = 0

$ws.Range("A28").Value = "TLEVCPQQHYCYDDHATSLYQPLFPQGPRMDINIWLWLSMPLNLHELRIWCAEDNGVWPHNSWKNPRKCNVVVTQPDTPPGS"
$ws.Range("B28").Value = 0
$ws.Range("C28").Value = 0

$ws.Range("A29").Value = "ITSVSWCTPGCTSEGGGSGCSHCC"
$ws.Range("B29").Value = 1
$ws.Range("C29").Value = 1

Write-Output "Updated rows 2-29 of training_set_prediction sheet"
